$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all_pathway_genes")

# Update the min/max length limits for each subunit length group
$ws.Range("G2:G5").Value = 1100
$ws.Range("H2:H5").Value = 1500

$ws.Range("G6:G9").Value = 500
$ws.Range("H6:H9").Value = 750

$ws.Range("G10:G13").Value = 300
$ws.Range("H10:H13").Value = 500

# Update the active selection to match the new review location
$ws.Activate()
$ws.Range("G11:H13").Select()
